$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "comments"
$ws.Range("C2").Value = "present"
$ws.Range("B1").Value = "AMOUNT"
$ws.Range("B2").Value = 500

$ws.Range("B3").Select()
